# ASK-41 Also read common sections again
#
# The row containing the "0" signalling-class entry (row 17: D17="0", E17=B4)
# is removed entirely from the "Normen en duidingsklassen" sheet, which
# shifts every row below it up by one. A few formulas in the rows above the
# deleted row are also updated to reflect new thresholds.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Normen en duidingsklassen")

# Update the "Signaleringskans per klasse" formulas (rows 13-15, unaffected by
# the row deletion below since they sit above row 17).
$ws.Range("E13").Formula = "=1/1000*B3"
$ws.Range("E14").Formula = "=1/100*B3"
$ws.Range("E15").Formula = "=1/10*B3"

# Remove the "0" row (old row 17) entirely; rows below shift up by one.
$ws.Rows.Item(17).Delete()

# The "-I" row (now row 17, was row 18) drops its "3*" multiplier.
$ws.Range("E17").Formula = "=B4"

# Reselect the cell the author ended up on before saving.
$ws.Range("E20").Select()
